# Week 15 logged + Week 16 simulated: append new per-play numbers to the
# running logs on YDS / ST, and update the season-to-date summary cells on
# OFF / DEF / ST / TURNS / PEN accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS — append newly logged run (row R) / pass (row P) play yardages
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value() + " 3 10 3 0 7 12 16 0 4 1 2 -1 8 8 5 1 1 9 3 -1 2 5 -1 10"
$ws.Range("C2").Value = $ws.Range("C2").Value() + " 2 3 6 2 0 3 4 2 0 8 3 6 3 4 4 3 6 3 -1 0 2 5 4 2 1 0"
$ws.Range("B3").Value = $ws.Range("B3").Value() + " 30 11 12 21 9 8 3 6 5 10 11 2 11 7 5 5 7 5 5 11"
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 10 5 6 3 7 22 7 7 1 7 6 12 11 33 2 6 18 3 43"

# ---------------------------------------------------------------------
# OFF — updated season totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = 135
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 45
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 21
$ws.Range("L2").Value = 272
$ws.Range("M2").Value = 162
$ws.Range("O2").Value = 27
$ws.Range("P2").Value = 15
$ws.Range("Q2").Value = 483

$ws.Range("B3").Value = 24
$ws.Range("C3").Value = 148
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 78
$ws.Range("G3").Value = 29
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 48
$ws.Range("J3").Value = 45
$ws.Range("N3").Value = 17

# ---------------------------------------------------------------------
# DEF — updated season totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 188
$ws.Range("F2").Value = 58
$ws.Range("G2").Value = 53
$ws.Range("J2").Value = 31
$ws.Range("L2").Value = 265
$ws.Range("M2").Value = 187
$ws.Range("Q2").Value = 501

$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 129
$ws.Range("E3").Value = 29
$ws.Range("F3").Value = 92
$ws.Range("G3").Value = 36
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 51
$ws.Range("J3").Value = 55
$ws.Range("N3").Value = 14

# ---------------------------------------------------------------------
# ST — updated season totals + appended kickoff/punt logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 44
$ws.Range("D2").Value = 62
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = 17
$ws.Range("J2").Value = 46
$ws.Range("K2").Value = 40

$ws.Range("B3").Value = 19
$ws.Range("D3").Value = $ws.Range("D3").Value() + " 54 40 54 40 49"

$ws.Range("B4").Value = $ws.Range("B4").Value() + " 63 58 55"
$ws.Range("D4").Value = $ws.Range("D4").Value() + " 16 3 0 0 14"

$ws.Range("B5").Value = $ws.Range("B5").Value() + " 98 0 23"
$ws.Range("D5").Value = $ws.Range("D5").Value() + " -2 48 3 0"

$ws.Range("B6").Value = $ws.Range("B6").Value() + " 20 16"

# ---------------------------------------------------------------------
# TURNS — updated season totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 6

# ---------------------------------------------------------------------
# PEN — updated season totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B4").Value = 3
